$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the formatting of the other header cells (copy from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add numeric values in the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
